# Add new columns I ("I0") and J ("IF") to the sheet, mirroring the
# header style already used by column H, and fill in the data rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Headers -----------------------------------------------------------
# Copy the existing header formatting (bold, bordered, centered) from H1
# onto the two new header cells, then set their text.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# --- Data ----------------------------------------------------------------
# I2:J73 values taken from the diff (row, I-value, J-value)
$data = @(
    @(6,7),
    @(8,8),
    @(8,8),
    @(8,8),
    @(7,7),
    @(9,9),
    @(8,8),
    @(6,7),
    @(12,12),
    @(8,8),
    @(8,8),
    @(5,6),
    @(9,9),
    @(9,9),
    @(7,7),
    @(8,8),
    @(10,10),
    @(9,9),
    @(8,8),
    @(9,9),
    @(9,9),
    @(9,9),
    @(7,7),
    @(7,7),
    @(8,8),
    @(8,8),
    @(9,9),
    @(8,8),
    @(8,8),
    @(8,8),
    @(8,8),
    @(8,8),
    @(9,9),
    @(8,8),
    @(9,9),
    @(9,9),
    @(7,7),
    @(9,9),
    @(7,7),
    @(8,8),
    @(9,9),
    @(8,8),
    @(9,9),
    @(9,9),
    @(8,8),
    @(8,9),
    @(8,9),
    @(9,9),
    @(9,9),
    @(9,9),
    @(8,8),
    @(9,9),
    @(8,9),
    @(8,8),
    @(9,9),
    @(9,9),
    @(7,7),
    @(7,7),
    @(9,9),
    @(9,9),
    @(8,8),
    @(8,9),
    @(9,9),
    @(8,8),
    @(5,5),
    @(7,7),
    @(6,6),
    @(6,7),
    @(5,5),
    @(5,5),
    @(3,3),
    @(3,3)
)

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 9).Value = $data[$i][0]
    $ws.Cells.Item($row, 10).Value = $data[$i][1]
}
